$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.856.65"
$ws.Range("D3").Value = "2.274.15"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'304.33"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").Value = "'93.15"
$ws.Range("E6").Value = "  +6.29%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").Value = "'32.63"
$ws.Range("E10").Value = "  +5.86%  "
$ws.Range("D11").Value = "'53.56"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "2.626.90"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "'14.28"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "2.284.59"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "'0.764"
$ws.Range("E18").Value = "  +3.62%  "
$ws.Range("D19").Value = "41.793.98"
$ws.Range("E20").Value = "  +8.55%  "
$ws.Range("D21").Value = "0.0₃0911"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("E22").Value = "  +3.12%  "
$ws.Range("D23").Value = "'67.38"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'243.70"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +4.73%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +5.38%  "
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").Value = "'9.63"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "'34.16"
$ws.Range("D32").Value = "'158.61"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("D35").Value = "'0.0753"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").Value = "'3.06"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").Value = "'16.70"
$ws.Range("E38").Value = "  +8.06%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "  +5.24%  "
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("D42").Value = "'3.94"
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("D43").Value = "2.072.02"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "'19.84"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +6.90%  "
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "'73.05"
$ws.Range("E50").Value = "  +7.58%  "
$ws.Range("E51").Value = "  +3.64%  "
